$wb = $excel.ActiveWorkbook

# Add the new "AddingToCart" worksheet after the last existing sheet (SearchFunction)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "AddingToCart"

# Header row, reusing the existing "Assertions" shared string
$newSheet.Range("A1").Value = "Assertions"

# New assertion text for the "adding to cart" test case
$newSheet.Range("A2").Value = "Product successfully added to your shopping cart"
$newSheet.Range("A2").WrapText = $true
$newSheet.Rows.Item(2).RowHeight = 87
$newSheet.Columns.Item(1).ColumnWidth = 8.7265625

$newSheet.PageSetup.Orientation = 1

$newSheet.Range("A2").Select() | Out-Null
